# The edit adds the per-region "additional spending" figures (65,000 USD each)
# for the five regions/zones that previously had blank entries in column B
# of the "Current expenditure" sheet, and highlights those cells in yellow
# to flag them as newly entered values. Finally the cursor is left on B14,
# matching the author's final selection when they saved the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("B5", "B6", "B9", "B10", "B14")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $cell.Value = 65000
    $cell.Interior.Color = 65535   # RGB(255,255,0) yellow, stored as BGR by COM
}

$ws.Range("B14").Select()
